$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D7").Value = -7.082100000000001
$ws.Range("C8").Value = -12.87749999999999
$ws.Range("C10").Value = -13.17849999999999
$ws.Range("C12").Value = -10.5506
$ws.Range("D14").Value = -7.798800000000004
$ws.Range("D15").Value = -7.965800000000001
$ws.Range("C18").Value = -13.67929999999999
$ws.Range("D18").Value = -8.85259999999999
$ws.Range("D20").Value = -7.552499999999998
$ws.Range("C25").Value = -13.50609999999999
$ws.Range("D29").Value = -7.171899999999999
$ws.Range("D30").Value = -7.265800000000002
$ws.Range("D31").Value = -8.453599999999996
$ws.Range("D35").Value = -8.216399999999991
$ws.Range("C37").Value = -13.7095
$ws.Range("D40").Value = -7.749199999999994
$ws.Range("D44").Value = -7.359499999999998
$ws.Range("D50").Value = -8.136999999999997
$ws.Range("D54").Value = -8.304900000000005
$ws.Range("C55").Value = -13.74509999999999
$ws.Range("C68").Value = -10.9785
$ws.Range("D68").Value = -6.991399999999996
$ws.Range("D76").Value = -7.312199999999998
$ws.Range("C77").Value = -12.4075
$ws.Range("C78").Value = -12.48520000000001
$ws.Range("C79").Value = -11.9553
$ws.Range("C80").Value = -13.2985
$ws.Range("C81").Value = -12.8754
$ws.Range("C82").Value = -12.59539999999999
$ws.Range("C84").Value = -13.39089999999999
$ws.Range("D87").Value = -8.100799999999994
$ws.Range("D88").Value = -7.452899999999997
$ws.Range("D92").Value = -6.3188
$ws.Range("D96").Value = -7.415300000000002
$ws.Range("D98").Value = -8.5395
$ws.Range("C101").Value = -13.1851
$ws.Range("D101").Value = -8.0426
$ws.Range("C102").Value = -13.72910000000001
$ws.Range("D102").Value = -7.6268
